$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @("Before", "Antes"),
    @("During", "Durante"),
    @("After", "Después"),
    @("What's Happened?", "¿Qué ha pasado?"),
    @("What's the Worst?", "¿Cuál es el peor?"),
    @("Cascadia Quake", "Terremoto de Cascadia"),
    @("Tsunami Zone", "Zona de tsunami"),
    @("If the dams failed", "Si fallaran las presas")
)

$row = 94
foreach ($pair in $pairs) {
    $ws.Range("A$row").Value = $pair[0]
    $ws.Range("B$row").Value = $pair[1]
    $row++
}
